$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The mapping rows for "RegistrationData.*" concepts shift down by one
# position (row 12 -> row 22, 22 -> 23, 23 -> 24, 24 -> 25) in columns
# B..J (cardinality/type/etc.), while column A and columns K..R of row 12
# stay untouched. Row 12 itself becomes blank in B..J, and a new row 25
# is appended carrying the data that used to live in row 24.

# --- Step 1: capture current contents of the rows that will move ---
$row12 = $ws.Range("B12:J12").Value2
$row22 = $ws.Range("B22:J22").Value2
$row23 = $ws.Range("B23:J23").Value2
$row24 = $ws.Range("B24:J24").Value2

# --- Step 2: clear row 12 (B..J) ---
$ws.Range("B12:J12").Value = ""

# --- Step 3: shift rows 22 -> 23 -> 24 -> 25, and old row 12 -> row 22 ---
$ws.Range("B25:J25").Value2 = $row24
$ws.Range("B24:J24").Value2 = $row23
$ws.Range("B23:J23").Value2 = $row22
$ws.Range("B22:J22").Value2 = $row12
